# Quarterly indexing esoteric bug-fix operation
#
# Column A holds quarterly dates that were incorrectly stored as the 1st
# day of the quarter-start month (Jan/Apr/Jul/Oct 1st). This fixes each
# date by shifting it forward to the 15th of the following month
# (e.g. 1988-07-01 -> 1988-08-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 150 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($serial -eq $null) { continue }

    $oldDate = [datetime]::FromOADate($serial)
    $newDate = $oldDate.AddMonths(1).AddDays(14)

    $cell.Value = $newDate
}
